$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: compareType changes from "<=" to "=", parameter from 20 to 0 ---
$ws.Range("E3").Value = "'="
$ws.Range("F3").Value = 0

# --- New row 6: canSignUp1 ---
$ws.Range("A6").Value = "canSignUp1"
$ws.Range("B6").Value = "是否可以签约条件1"
$ws.Range("C6").Value = "city"
$ws.Range("D6").Value = "totalPercentage"
$ws.Range("E6").Value = "<"
$ws.Range("F6").Value = 100

# --- New row 7: canSignUp2 ---
$ws.Range("A7").Value = "canSignUp2"
$ws.Range("B7").Value = "是否可以签约条件2"
$ws.Range("C7").Value = "city"
$ws.Range("D7").Value = "guildNumber"
$ws.Range("E7").Value = "<"
$ws.Range("F7").Value = 3

# --- New row 8: canSignUp ---
$ws.Range("A8").Value = "canSignUp"
$ws.Range("B8").Value = "是否可以签约"
$ws.Range("C8").Value = "and"
$ws.Range("D8").Value = ";"
$ws.Range("E8").Value = ";"
$ws.Range("F8").Value = "canSignUp1;canSignUp2"

# --- New row 9: canSignUpMoneyEnough ---
$ws.Range("A9").Value = "canSignUpMoneyEnough"
$ws.Range("B9").Value = "签约金够"
$ws.Range("C9").Value = "city"
$ws.Range("D9").Value = "signUpMoney"
$ws.Range("E9").Value = "<="
$ws.Range("F9").Value = "money"

# --- New row 10: canMilitaryInvestMoneyEnough ---
$ws.Range("A10").Value = "canMilitaryInvestMoneyEnough"
$ws.Range("B10").Value = "军事投资钱是否够"
$ws.Range("C10").Value = "city"
$ws.Range("D10").Value = "militaryInvestMoney"
$ws.Range("E10").Value = "<="
$ws.Range("F10").Value = "money"

# --- Apply the "Chinese" font style (matches existing condition-id column cells) ---
$ws.Range("A6").Font.Name = "宋体"
$ws.Range("A7").Font.Name = "宋体"
$ws.Range("A8").Font.Name = "宋体"
$ws.Range("A9").Font.Name = "宋体"
$ws.Range("F8").Font.Name = "宋体"

# --- Column widths (best-fit-ish, mirrors author resizing columns A, B and adding F) ---
# Values chosen to land in the middle of this engine's column-width quantization
# bucket that's closest to the authored widths (27, 18.1640625, 22.1640625 chars).
$ws.Columns.Item(1).ColumnWidth = 26.17
$ws.Columns.Item(2).ColumnWidth = 17.33
$ws.Columns.Item(6).ColumnWidth = 21.33

# --- Selection moves to B9 ---
$ws.Range("B9").Select()
